# "Semic Automatic Tool added" - extend the dbpedia extract sheet with
# birthDate / birthPlace / deathDate columns (D:F) for the first two
# people, drop the now-unused 4th data row (Y._D._Tiwari), and shrink
# the used range down to A1:F3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells (row 1) -------------------------------------------
$ws.Range("D1").Value = "http://dbpedia.org/ontology/birthDate"
$ws.Range("E1").Value = "http://dbpedia.org/ontology/birthPlace"
$ws.Range("F1").Value = "http://dbpedia.org/ontology/deathDate"

# Match the header formatting already used for A1:C1 (bold, bordered,
# centered/top-aligned) by copying it across onto the new header cells.
$ws.Range("A1").Copy()
$ws.Range("D1:F1").PasteSpecial(-4122)

# --- Row 2 (Reginald Pole) new data --------------------------------------
$ws.Range("D2").Value = "1500-03-12 "
$ws.Range("E2").Value = "http://dbpedia.org/resource/Kingdom_of_England http://dbpedia.org/resource/Staffordshire http://dbpedia.org/resource/Stourton_Castle "
$ws.Range("F2").Value = "1558-11-17 "

# --- Row 3 (R. R. Sundara Rao) new data ----------------------------------
# D3 and F3 stay blank for this record but remain present cells.
$ws.Range("D3").Value = " "
$ws.Range("E3").Value = "http://dbpedia.org/resource/Andhra_Pradesh "
$ws.Range("F3").Value = " "

# --- Drop row 4 (Y._D._Tiwari) entirely ----------------------------------
$ws.Rows.Item(4).Delete()
